$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 197 (PAP / Port-au-Prince, Haiti) - all following rows shift up by one.
$ws.Rows.Item(197).Delete()
